$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data block for rows 2-5 with the data block for rows 6-9
# (columns D, M, N, O, P, R, S), effectively re-ordering the two weekly
# batches of prices while keeping the quality grade (column L) aligned to
# its row.

$rowsTop = 2,3,4,5
$rowsBottom = 6,7,8,9

for ($i = 0; $i -lt 4; $i++) {
    $rTop = $rowsTop[$i]
    $rBot = $rowsBottom[$i]

    # Capture current ("before") values for both rows.
    # NOTE: use Value2 for reads - the Value getter on this runtime returns
    # reflection metadata instead of the actual cell content.
    $dTop = $ws.Range("D$rTop").Value2
    $mTop = $ws.Range("M$rTop").Value2
    $nTop = $ws.Range("N$rTop").Value2
    $oTop = $ws.Range("O$rTop").Value2
    $pTop = $ws.Range("P$rTop").Value2
    $rTopOrigin = $ws.Range("R$rTop").Value2
    $sTop = $ws.Range("S$rTop").Value2

    $dBot = $ws.Range("D$rBot").Value2
    $mBot = $ws.Range("M$rBot").Value2
    $nBot = $ws.Range("N$rBot").Value2
    $oBot = $ws.Range("O$rBot").Value2
    $pBot = $ws.Range("P$rBot").Value2
    $rBotOrigin = $ws.Range("R$rBot").Value2
    $sBot = $ws.Range("S$rBot").Value2

    # Write bottom's original values into top row
    $ws.Range("D$rTop").Value = $dBot
    $ws.Range("M$rTop").Value = $mBot
    $ws.Range("N$rTop").Value = $nBot
    $ws.Range("O$rTop").Value = $oBot
    $ws.Range("P$rTop").Value = $pBot
    $ws.Range("R$rTop").Value = $rBotOrigin
    $ws.Range("S$rTop").Value = $sBot

    # Write top's original values into bottom row
    $ws.Range("D$rBot").Value = $dTop
    $ws.Range("M$rBot").Value = $mTop
    $ws.Range("N$rBot").Value = $nTop
    $ws.Range("O$rBot").Value = $oTop
    $ws.Range("P$rBot").Value = $pTop
    $ws.Range("R$rBot").Value = $rTopOrigin
    $ws.Range("S$rBot").Value = $sTop
}
